$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.23 = 8511.14 pesos`n✅ 8511.14 pesos = 2.21 = 941.01 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 449
$wsTasas.Range("O10").Value = 3821.5

$wsTasas.Range("N12").Value = 3844.9
$wsTasas.Range("O12").Value = 425.1
